$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Ensure the GDP value column stores these numeric-looking figures as text,
# consistent with the existing data (E2:E60 already hold text values).
$ws.Range("E2:E68").NumberFormat = "@"

# Updated GDP per Capita figures for existing years (1950-2008)
$gdpValues = @{
    2 = "1307"
    3 = "1339"
    4 = "1369"
    5 = "1400"
    6 = "1431"
    7 = "1462"
    8 = "1492"
    9 = "1522"
    10 = "1553"
    11 = "1581"
    12 = "1567"
    13 = "1597"
    14 = "1525"
    15 = "1492"
    16 = "1486"
    17 = "1466"
    18 = "1427"
    19 = "1452"
    20 = "1435"
    21 = "1476"
    22 = "1508"
    23 = "1427"
    24 = "1368"
    25 = "1403"
    26 = "1333"
    27 = "1323"
    28 = "1345"
    29 = "1417"
    30 = "1398"
    31 = "1326"
    32 = "1234"
    33 = "1219"
    34 = "1208"
    35 = "1105"
    36 = "1146"
    37 = "1148"
    38 = "1162"
    39 = "1079"
    40 = "1073"
    41 = "1070"
    42 = "1023"
    43 = "1006.99580160323"
    44 = "965.193376349712"
    45 = "933.060174663238"
    46 = "952.017417741463"
    47 = "978.556614747978"
    48 = "885.886400596378"
    49 = "917.47997639614"
    50 = "936.238262138212"
    51 = "967.843203629496"
    52 = "963.946180555058"
    53 = "952.511809067715"
    54 = "928.771638513173"
    55 = "845.505545088412"
    56 = "853.182021590341"
    57 = "858.523700017484"
    58 = "881.485165116962"
    59 = "897.478307013892"
    60 = "901.851345375175"
}
foreach ($row in $gdpValues.Keys) {
    $ws.Cells.Item($row, 5).Value = $gdpValues[$row]
}

# Append new yearly rows for 2009-2016
$newRows = @(
    @{Row=61; Year=2009; Value="898.700480180656"}
    @{Row=62; Year=2010; Value="910.310378860757"}
    @{Row=63; Year=2011; Value="909"}
    @{Row=64; Year=2012; Value="916"}
    @{Row=65; Year=2013; Value="567"}
    @{Row=66; Year=2014; Value="561"}
    @{Row=67; Year=2015; Value="576"}
    @{Row=68; Year=2016; Value="589"}
)
foreach ($item in $newRows) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = 140
    $ws.Cells.Item($r, 2).Value = "Central African Republic"
    $ws.Cells.Item($r, 3).Value = "GDP per Capita"
    $ws.Cells.Item($r, 4).Value = $item.Year
    $ws.Cells.Item($r, 5).Value = $item.Value
}